# Weekly update: add a new week's worth of "Acelga" price data (Vega
# Monumental Concepción) at the front of the date-sorted block that starts
# at row 386. This pushes the existing rows 386:440 down to 388:442 and
# extends the used range from A1:R440 to A1:R442.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 386 - everything currently at row 386 and below
# shifts down by two rows, matching rows 388:442 in the target file.
$ws.Rows("386:387").Insert()

# New row 386: "Primera" quality entry for the new week (2023-07-20, serial 45127)
$ws.Cells.Item(386, 1).Value = 11
$ws.Cells.Item(386, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(386, 3).Value = "Bíobío"
$ws.Cells.Item(386, 4).Value = 45127
$ws.Cells.Item(386, 5).Value = 8
$ws.Cells.Item(386, 6).Value = 100112009
$ws.Cells.Item(386, 7).Value = "Acelga"
$ws.Cells.Item(386, 8).Value = "Sin especificar"
$ws.Cells.Item(386, 9).Value = "Primera"
$ws.Cells.Item(386, 10).Value = 100
$ws.Cells.Item(386, 11).Value = 600
$ws.Cells.Item(386, 12).Value = 600
$ws.Cells.Item(386, 13).Value = 600
$ws.Cells.Item(386, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(386, 15).Value = "Región de Ñuble"
$ws.Cells.Item(386, 16).Value = 600
$ws.Cells.Item(386, 17).Value = 1
$ws.Cells.Item(386, 18).Value = "Hortaliza"

# New row 387: "Segunda" quality entry for the same new week
$ws.Cells.Item(387, 1).Value = 11
$ws.Cells.Item(387, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(387, 3).Value = "Bíobío"
$ws.Cells.Item(387, 4).Value = 45127
$ws.Cells.Item(387, 5).Value = 8
$ws.Cells.Item(387, 6).Value = 100112009
$ws.Cells.Item(387, 7).Value = "Acelga"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Segunda"
$ws.Cells.Item(387, 10).Value = 80
$ws.Cells.Item(387, 11).Value = 650
$ws.Cells.Item(387, 12).Value = 650
$ws.Cells.Item(387, 13).Value = 650
$ws.Cells.Item(387, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(387, 15).Value = "Región de Ñuble"
$ws.Cells.Item(387, 16).Value = 650
$ws.Cells.Item(387, 17).Value = 1
$ws.Cells.Item(387, 18).Value = "Hortaliza"

# Keep the date columns formatted the same as the rest of column D.
$ws.Range("D386:D387").NumberFormat = $ws.Range("D388").NumberFormat
